$wb = $excel.ActiveWorkbook

# Sheet "展览" updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 107
$ws1.Range("F4").Value = 132
$ws1.Range("F5").Value = 2868
$ws1.Range("F6").Value = 287
$ws1.Range("F7").Value = 392

# Sheet "全部类型" updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 107
$ws4.Range("F4").Value = 132
$ws4.Range("F5").Value = 2868
$ws4.Range("F6").Value = 287
$ws4.Range("F9").Value = 392
